$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$wsOut = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4245-RBI-EI-DB-SAR-REC-RNI-FEE-FFConMONonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-ONTIME-PE-1st"

# Update the product name (B1) - rename part of the test name.
# Update it on both sheets (they share the same underlying string) so the
# shared string table keeps a single shared entry instead of duplicating it.
$ws.Range("B1").Value = $newProductName
$wsOut.Range("B1").Value = $newProductName

# Update the short name (B2) from numeric productid 4245 to the new string id "424a"
# (removes the inter-test dependency on a hard-coded numeric product id)
$ws.Range("B2").Value = "424a"

# Move the active selection to B7 on the input sheet, then restore the
# output sheet as the active tab (matches the workbook's original active tab).
$ws.Range("B7").Select()
$wsOut.Activate()
